$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected, which blocks direct cell writes, so it must
# be unprotected before the updated figures below can be entered.
$ws.Unprotect()

# Update the confidential disclaimer date from 2021-03-17 to 2021-03-18,
# preserving the embedded line break between the two sentences.
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + "`n" + "Model holdings provided as of 2021-03-18 for illustrative purposes only and are subject to change."

# Update the recalculated Weight (column D) and Percent Change (column E)
# figures for the three holding rows.
$ws.Range("D2").Value = 0.8412809519643013
$ws.Range("E2").Value = -0.009100364014560602

$ws.Range("D3").Value = 0.1587190480356987
$ws.Range("E3").Value = -0.01837559720690918

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = -0.01057252019715926

$wb.Save()
